$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 16
$ws.Range("I12").Value = 0
$ws.Range("J12").Value = 16
$ws.Range("K12").Value = 0
$ws.Range("L12").ClearContents()
$ws.Range("M12").Value = 16
$ws.Range("N12").Value = -356
$ws.Range("H55").Value = 841.5833
$ws.Range("I55").Value = 401.375
$ws.Range("K55").Value = 401.375
$ws.Range("M55").Value = -187.375
$ws.Range("H70").Value = 4339
$ws.Range("I70").Value = 1684
$ws.Range("J70").Value = 5856.143
$ws.Range("K70").Value = 5052
$ws.Range("L70").Value = 17568.429
$ws.Range("M70").Value = -4782
$ws.Range("N70").Value = -18108.429
$ws.Range("H73").Value = 4339
$ws.Range("I73").Value = 1684
$ws.Range("J73").Value = 5856.143
$ws.Range("K73").Value = 5052
$ws.Range("L73").Value = 17568.429
$ws.Range("M73").Value = -4116
$ws.Range("N73").Value = -19440.429
$ws.Range("H80").Value = 583.3333
$ws.Range("I80").Value = 583.3333
$ws.Range("J80").Value = 0
$ws.Range("K80").Value = 1749.9999
$ws.Range("L80").Value = 0
$ws.Range("M80").ClearContents()
$ws.Range("N80").Value = -751.9999
$ws.Range("H83").Value = 583.3333
$ws.Range("I83").Value = 583.3333
$ws.Range("J83").Value = 0
$ws.Range("K83").Value = 5249.9997
$ws.Range("L83").Value = 0
$ws.Range("M83").ClearContents()
$ws.Range("N83").Value = -257.9997000000003
$ws.Range("H86").Value = 5665
$ws.Range("I86").Value = 4544.8184
$ws.Range("K86").Value = 4544.8184
$ws.Range("M86").Value = -3421.8184
$ws.Range("H89").Value = 5665
$ws.Range("I89").Value = 4544.8184
$ws.Range("K89").Value = 22724.092
$ws.Range("M89").Value = -17108.092
$ws.Range("H99").Value = 4604.077
$ws.Range("I99").Value = 2238
$ws.Range("J99").Value = 5313.9
$ws.Range("K99").Value = 6714
$ws.Range("L99").Value = 15941.7
$ws.Range("M99").Value = -5216
$ws.Range("N99").Value = -18937.7
$ws.Range("H100").Value = 2965.8333
$ws.Range("J100").Value = 1745
$ws.Range("L100").Value = 1745
$ws.Range("N100").Value = -2827
$ws.Range("H111").Value = 1524
$ws.Range("I111").Value = 383
$ws.Range("J111").Value = 1850
$ws.Range("K111").Value = 1149
$ws.Range("L111").Value = 5550
$ws.Range("M111").Value = 1918
$ws.Range("N111").Value = -11684
$ws.Range("H112").Value = 2019.5
$ws.Range("I112").Value = 2500
$ws.Range("K112").Value = 7500
$ws.Range("M112").Value = -6392
$ws.Range("H113").Value = 7898.3335
$ws.Range("I113").Value = 3705
$ws.Range("K113").Value = 3705
$ws.Range("M113").Value = -451
$ws.Range("H129").Value = 1428.6
$ws.Range("I129").Value = 1035.75
$ws.Range("K129").Value = 3107.25
$ws.Range("M129").Value = 1892.75
$ws.Range("H135").Value = 900
$ws.Range("I135").Value = 900
$ws.Range("K135").Value = 8100
$ws.Range("M135").Value = -5565

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H50").Value = 1300
$ws.Range("I50").Value = 600
$ws.Range("K50").Value = 600
$ws.Range("M50").Value = 114
$ws.Range("H61").Value = 3736.2173
$ws.Range("I61").Value = 2946.65
$ws.Range("K61").Value = 2946.65
$ws.Range("M61").Value = -2734.65
$ws.Range("H101").Value = 59997.5
$ws.Range("J101").Value = 59997.5
$ws.Range("L101").Value = 59997.5
$ws.Range("N101").Value = -66487.5
$ws.Range("H136").Value = 3736.2173
$ws.Range("I136").Value = 2946.65
$ws.Range("K136").Value = 8839.950000000001
$ws.Range("M136").Value = -6289.950000000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H25").Value = 2304.6667
$ws.Range("I25").Value = 2304.6667
$ws.Range("K25").Value = 2304.6667
$ws.Range("M25").Value = -2069.6667
$ws.Range("H48").Value = 240000
$ws.Range("J48").Value = 240000
$ws.Range("L48").Value = 240000
$ws.Range("N48").Value = -240830
$ws.Range("H86").Value = 4054.1
$ws.Range("I86").Value = 1889.2727
$ws.Range("J86").Value = 6700
$ws.Range("K86").Value = 1889.2727
$ws.Range("L86").Value = 6700
$ws.Range("M86").Value = -766.2727
$ws.Range("N86").Value = -8946
$ws.Range("H89").Value = 4054.1
$ws.Range("I89").Value = 1889.2727
$ws.Range("J89").Value = 6700
$ws.Range("K89").Value = 9446.363499999999
$ws.Range("L89").Value = 33500
$ws.Range("M89").Value = -3830.363499999999
$ws.Range("N89").Value = -44732
$ws.Range("H134").Value = 3324.6667
$ws.Range("I134").Value = 3000
$ws.Range("K134").Value = 9000
$ws.Range("M134").Value = -6465

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 4688.385
$ws.Range("I58").Value = 4222.6665
$ws.Range("K58").Value = 4222.6665
$ws.Range("M58").Value = -4019.6665
$ws.Range("H62").Value = 3500
$ws.Range("I62").Value = 3500
$ws.Range("K62").Value = 3500
$ws.Range("M62").Value = -2876
$ws.Range("H65").Value = 3500
$ws.Range("I65").Value = 3500
$ws.Range("K65").Value = 17500
$ws.Range("M65").Value = -14380
$ws.Range("H74").Value = 45917.223
$ws.Range("J74").Value = 47906.875
$ws.Range("L74").Value = 47906.875
$ws.Range("N74").Value = -49654.875
$ws.Range("H77").Value = 45917.223
$ws.Range("J77").Value = 47906.875
$ws.Range("L77").Value = 143720.625
$ws.Range("N77").Value = -152456.625
$ws.Range("H105").Value = 4000.5
$ws.Range("J105").Value = 3011
$ws.Range("L105").Value = 3011
$ws.Range("N105").Value = -6505
$ws.Range("H136").Value = 4688.385
$ws.Range("I136").Value = 4222.6665
$ws.Range("K136").Value = 12667.9995
$ws.Range("M136").Value = -10117.9995

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 47.77778
$ws.Range("I2").Value = 16.666666
$ws.Range("K2").Value = 99.999996
$ws.Range("M2").Value = 13.000004
$ws.Range("H7").Value = 262.5
$ws.Range("J7").Value = 262.5
$ws.Range("L7").Value = 787.5
$ws.Range("N7").Value = -1011.5
$ws.Range("H22").Value = 295
$ws.Range("J22").Value = 295
$ws.Range("L22").Value = 885
$ws.Range("N22").Value = -1223
$ws.Range("H23").Value = 500
$ws.Range("J23").Value = 500
$ws.Range("L23").Value = 1500
$ws.Range("N23").Value = -1970
$ws.Range("H27").Value = 295
$ws.Range("J27").Value = 295
$ws.Range("L27").Value = 885
$ws.Range("N27").Value = -1089
$ws.Range("H54").Value = 1350
$ws.Range("I54").Value = 700
$ws.Range("K54").Value = 2100
$ws.Range("M54").Value = -1541
$ws.Range("H129").Value = 1199
$ws.Range("I129").Value = 0
$ws.Range("K129").Value = 0
$ws.Range("M129").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 212.36363
$ws.Range("J2").Value = 289
$ws.Range("L2").Value = 289
$ws.Range("N2").Value = -515
$ws.Range("H80").Value = 3065
$ws.Range("J80").Value = 2966.6667
$ws.Range("L80").Value = 2966.6667
$ws.Range("N80").Value = -4962.6667
$ws.Range("H83").Value = 3065
$ws.Range("J83").Value = 2966.6667
$ws.Range("L83").Value = 14833.3335
$ws.Range("N83").Value = -24817.3335
$ws.Range("H132").Value = 2284.4285
$ws.Range("I132").Value = 2415.1667
$ws.Range("J132").Value = 1500
$ws.Range("K132").Value = 7245.500100000001
$ws.Range("L132").Value = 4500
$ws.Range("M132").Value = -4715.500100000001
$ws.Range("N132").Value = -9560

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("M22").ClearContents()
$ws.Range("H27").Value = 0
$ws.Range("I27").Value = 0
$ws.Range("K27").Value = 0
$ws.Range("M27").ClearContents()
$ws.Range("H55").Value = 805.5
$ws.Range("I55").Value = 805.5
$ws.Range("J55").Value = 0
$ws.Range("K55").Value = 805.5
$ws.Range("L55").Value = 0
$ws.Range("M55").ClearContents()
$ws.Range("N55").Value = -632.5
$ws.Range("H136").Value = 2796
$ws.Range("I136").Value = 2796
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 8388
$ws.Range("L136").Value = 0
$ws.Range("M136").ClearContents()
$ws.Range("N136").Value = -5838

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 5333.3335
$ws.Range("I81").Value = 3000
$ws.Range("J81").Value = 10000
$ws.Range("K81").Value = 6000
$ws.Range("L81").Value = 20000
$ws.Range("M81").Value = -4939
$ws.Range("N81").Value = -22122
$ws.Range("H84").Value = 5333.3335
$ws.Range("I84").Value = 3000
$ws.Range("J84").Value = 10000
$ws.Range("K84").Value = 30000
$ws.Range("L84").Value = 100000
$ws.Range("M84").Value = -24696
$ws.Range("N84").Value = -110608
$ws.Range("H107").Value = 1756.5714
$ws.Range("I107").Value = 1972.6666
$ws.Range("J107").Value = 460
$ws.Range("K107").Value = 5917.9998
$ws.Range("L107").Value = 1380
$ws.Range("M107").Value = -3997.9998
$ws.Range("N107").Value = -5220
$ws.Range("H136").Value = 3765.1738
$ws.Range("I136").Value = 2407.4167
$ws.Range("K136").Value = 7222.250100000001
$ws.Range("M136").Value = -4672.250100000001
